# 3.18: Complete Advanced Algebra
# The deck starts with zero slides (no <p:sldIdLst>). This change inserts
# the first slide, using the "Title Slide" layout (ctrTitle + subTitle
# placeholders), leaving both placeholders empty for now.

$p = $ppt.ActivePresentation

# PpSlideLayout.ppLayoutTitle = 1 -> slideLayout1.xml ("标题幻灯片" / Title Slide),
# which is the layout that exposes the ctrTitle + subTitle placeholders.
$s = $p.Slides.Add(1, 1)

Write-Host "Added slide 1 (Title Slide layout); slide count = $($p.Slides.Count)"
